# database.xlsx ("Stickers" sheet) update:
#  - A2/B2 contents are swapped: A2 now holds "привет", B2 now holds the
#    (new) Telegram file_id-looking string.
#  - The new file_id string in B2 begins with a capital letter, but Excel
#    still stores the cell with an explicit quote-prefix ("Text that looks
#    like a formula/prefixed" flag -> <xf quotePrefix="1"/>), exactly like
#    the original B2/A2 content did. Typing a leading "'" via COM Value
#    reproduces that flag without the apostrophe becoming part of the text.
#  - A new, empty row 3 (B3) is added, also carrying the quote-prefix style.
#  - Column B is widened slightly and B2 becomes the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stickers")

# A2: was the file_id-looking text, now the greeting string.
$ws.Range("A2").Value = "привет"

# B2: was the greeting string, now the new file_id-looking text. The leading
# "'" forces Excel's text/quote-prefix handling (mirrors the source cell's
# original formatting) without becoming part of the stored string.
$ws.Range("B2").Value = "'CAACAgIAAxkBAAP6Yh3gGS64cUoFmCetDHxtn3oZ6VcAAkgBAAJ7TioQQ77E0Xv6lP8jBA"

# B3: brand-new, empty cell that still carries the quote-prefix style used
# above. Entering then clearing a lone "'" leaves the style flag applied to
# an otherwise value-less cell.
$ws.Range("B3").Value = "'"
$ws.Range("B3").Value = ""

# Column B grows a bit to comfortably fit the longer file_id string.
$ws.Columns("B").ColumnWidth = 89.3

# B2 becomes the selected / active cell on the sheet.
$ws.Range("B2").Select()
